# Update the four-row header of the decomposition table (row 4) to match
# the revised draft: the superscripts on P are capitalized (i->I, f->F,
# c->C), the "Val-Loan Dif" label is reworded to "Value-Loan", and the
# 1(Def_i) indicator cell is swapped to the column that used to hold the
# value/loan-difference label.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = '$\sum_t P^I_{it}$'
$ws.Range("D4").Value = '$\sum_t P^F_{it}$'
$ws.Range("E4").Value = '$\mathds{1}(\text{Def}_i)}\times\sum_t P^C_{it}$'
$ws.Range("F4").Value = '$\mathds{1}(\text{Def}_i)}\times \text{Value-Loan}_i$'
$ws.Range("G4").Value = '$\mathds{1}(\text{Def}_i)$'
